$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Category" in A1, matching the style of the other header cells (B1:W1)
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Remove the old header-style formatting from the data rows in column A (A2:A46),
# leaving them as plain/default-styled cells.
$ws.Range("A2:A46").ClearFormats()
